$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new paper entry (row 14): Number, Name, Path
$ws.Range("C14").Value = 11
$ws.Range("D14").Value = "11_Improving Circuit Performance in a Trapped-Ion"
$ws.Range("E14").Value = "C:\Jeonghyun\GIT\QSCOUT"

# Update selection to match the new active cell
$ws.Range("E14").Select()
